$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record this week's (Week 4) activity log entry for Cedric Stephani (column C, row 4)
$ws.Range("C4").Value = "Downloaded the systems programs like Visual Code, ReactJS/NodeJS, and MongoDB. Worked on user stories and project scheduling as well as started working on a rough database model which reflects systems/user stories needs."
